$d = $word.ActiveDocument

$replacements = @(
    @{old = "628×6="; new = "320×3="},
    @{old = "637×2="; new = "461×4="},
    @{old = "431×3="; new = "201×7="},
    @{old = "235×2="; new = "808×4="},
    @{old = "979×3="; new = "183×2="},
    @{old = "778×3="; new = "890×5="},
    @{old = "781×7="; new = "548×3="},
    @{old = "653×8="; new = "909×6="},
    @{old = "595×8="; new = "586×2="},
    @{old = "950×2="; new = "507×5="},
    @{old = "432×8="; new = "653×7="},
    @{old = "186×8="; new = "324×8="},
    @{old = "855×7="; new = "196×9="},
    @{old = "399×8="; new = "853×4="},
    @{old = "342×5="; new = "492×7="},
    @{old = "806×7="; new = "923×9="},
    @{old = "591×6="; new = "918×4="},
    @{old = "183×3="; new = "631×7="},
    @{old = "344×5="; new = "586×7="},
    @{old = "668×4="; new = "872×6="},
    @{old = "173×6="; new = "985×8="},
    @{old = "377×4="; new = "661×9="},
    @{old = "859×9="; new = "407×6="},
    @{old = "925×3="; new = "219×3="},
    @{old = "964×7="; new = "983×2="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}

$d.Save()
